$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add "NA" values under duplicate_image_filename (column E) for rows 2-21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
